# Commit: trial, block_num --> iTrial, iBlock
# Rename the two header labels on row 1 of Sheet1 (columns BB and BC) that
# previously read "trial" and "block_num" to "iTrial" and "iBlock".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("BB1").Value = "iTrial"
$ws.Range("BC1").Value = "iBlock"

# Mirror the author's resulting view state: scrolled right toward the
# renamed columns, with the description cell BC2 selected.
$excel.ActiveWindow.ScrollColumn = $ws.Range("AY1").Column
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("BC2").Select()
